$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4932.1665
$ws.Range("J19").Value = 6399.5
$ws.Range("L19").Value = 6399.5
$ws.Range("N19").Value = -6749.5
$ws.Range("H32").Value = 5106.625
$ws.Range("I32").Value = 2499
$ws.Range("J32").Value = 5479.143
$ws.Range("K32").Value = 2499
$ws.Range("L32").Value = 5479.143
$ws.Range("M32").Value = -2173
$ws.Range("N32").Value = -6131.143
$ws.Range("H40").Value = 3022.6191
$ws.Range("I40").Value = 2762.125
$ws.Range("J40").Value = 3856.2
$ws.Range("K40").Value = 2762.125
$ws.Range("L40").Value = 3856.2
$ws.Range("M40").Value = -2587.125
$ws.Range("N40").Value = -4206.2
$ws.Range("H132").Value = 6655.591
$ws.Range("I132").Value = 8073.5884
$ws.Range("J132").Value = 1834.4
$ws.Range("K132").Value = 24220.7652
$ws.Range("L132").Value = 5503.200000000001
$ws.Range("M132").Value = -21690.7652
$ws.Range("N132").Value = -10563.2
$ws.Range("H138").Value = 25002198
$ws.Range("I138").Value = 1126
$ws.Range("J138").Value = 58827180
$ws.Range("K138").Value = 3378
$ws.Range("L138").Value = 176481540
$ws.Range("M138").Value = 1762
$ws.Range("N138").Value = -176491820

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4070.157
$ws.Range("I32").Value = 4012.1458
$ws.Range("K32").Value = 4012.1458
$ws.Range("M32").Value = -3725.1458
$ws.Range("H74").Value = 6558.75
$ws.Range("I74").Value = 6583.0435
$ws.Range("K74").Value = 6583.0435
$ws.Range("M74").Value = -5709.0435
$ws.Range("H77").Value = 6558.75
$ws.Range("I77").Value = 6583.0435
$ws.Range("K77").Value = 32915.2175
$ws.Range("M77").Value = -28547.2175
$ws.Range("H102").Value = 4483.273
$ws.Range("I102").Value = 4007.1052
$ws.Range("K102").Value = 4007.1052
$ws.Range("M102").Value = -2385.1052
$ws.Range("H110").Value = 1120
$ws.Range("I110").Value = 1014.3333
$ws.Range("J110").Value = 1357.75
$ws.Range("K110").Value = 1014.3333
$ws.Range("L110").Value = 1357.75
$ws.Range("M110").Value = 1030.6667
$ws.Range("N110").Value = -5447.75
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null
$ws.Range("H114").Value = 74566.336
$ws.Range("J114").Value = 74566.336
$ws.Range("L114").Value = 74566.336
$ws.Range("N114").Value = -83244.336
$ws.Range("H132").Value = 32533.564
$ws.Range("I132").Value = 2335.5264
$ws.Range("J132").Value = 175974.25
$ws.Range("K132").Value = 7006.5792
$ws.Range("L132").Value = 527922.75
$ws.Range("M132").Value = -4476.5792
$ws.Range("N132").Value = -532982.75
$ws.Range("H140").Value = 78606.75
$ws.Range("J140").Value = 78606.75
$ws.Range("L140").Value = 78606.75
$ws.Range("N140").Value = -88966.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 5500
$ws.Range("J30").Value = 5500
$ws.Range("L30").Value = 5500
$ws.Range("N30").Value = -5750
$ws.Range("H86").Value = 9093.821
$ws.Range("I86").Value = 1819.1364
$ws.Range("K86").Value = 1819.1364
$ws.Range("M86").Value = -696.1364000000001
$ws.Range("H89").Value = 9093.821
$ws.Range("I89").Value = 1819.1364
$ws.Range("K89").Value = 9095.682000000001
$ws.Range("M89").Value = -3479.682000000001
$ws.Range("H105").Value = 17002.715
$ws.Range("I105").Value = 18169.834
$ws.Range("K105").Value = 18169.834
$ws.Range("M105").Value = -16422.834
$ws.Range("H134").Value = 1093.1666
$ws.Range("I134").Value = 1083.4546
$ws.Range("K134").Value = 3250.3638
$ws.Range("M134").Value = -715.3638000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 192.14285
$ws.Range("J7").Value = 254.10527
$ws.Range("L7").Value = 254.10527
$ws.Range("N7").Value = -480.10527
$ws.Range("H16").Value = 1540.7142
$ws.Range("J16").Value = 1995
$ws.Range("L16").Value = 1995
$ws.Range("N16").Value = -2569
$ws.Range("H31").Value = 4127.1836
$ws.Range("I31").Value = 2657.5417
$ws.Range("J31").Value = 5538.04
$ws.Range("K31").Value = 2657.5417
$ws.Range("L31").Value = 5538.04
$ws.Range("M31").Value = -2362.5417
$ws.Range("N31").Value = -6128.04
$ws.Range("H34").Value = 4127.1836
$ws.Range("I34").Value = 2657.5417
$ws.Range("J34").Value = 5538.04
$ws.Range("K34").Value = 2657.5417
$ws.Range("L34").Value = 5538.04
$ws.Range("M34").Value = -2455.5417
$ws.Range("N34").Value = -5942.04
$ws.Range("H105").Value = 11756.272
$ws.Range("I105").Value = 2700.889
$ws.Range("J105").Value = 52505.5
$ws.Range("K105").Value = 2700.889
$ws.Range("L105").Value = 52505.5
$ws.Range("M105").Value = -953.8890000000001
$ws.Range("N105").Value = -55999.5
$ws.Range("H113").Value = 1540.7142
$ws.Range("J113").Value = 1995
$ws.Range("L113").Value = 1995
$ws.Range("N113").Value = -6335
$ws.Range("H122").Value = 2485.7
$ws.Range("J122").Value = 3299.6667
$ws.Range("L122").Value = 9899.000100000001
$ws.Range("N122").Value = -14799.0001
$ws.Range("H132").Value = 3444.5918
$ws.Range("I132").Value = 3246.2856
$ws.Range("K132").Value = 9738.856800000001
$ws.Range("M132").Value = -7208.856800000001
$ws.Range("H134").Value = 4184.364
$ws.Range("I134").Value = 3624.2
$ws.Range("J134").Value = 4651.1665
$ws.Range("K134").Value = 10872.6
$ws.Range("L134").Value = 13953.4995
$ws.Range("M134").Value = -8337.599999999999
$ws.Range("N134").Value = -19023.4995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.47059
$ws.Range("I2").Value = 47.6
$ws.Range("J2").Value = 54.57143
$ws.Range("K2").Value = 285.6
$ws.Range("L2").Value = 327.42858
$ws.Range("M2").Value = -172.6
$ws.Range("N2").Value = -553.42858
$ws.Range("H23").Value = 691.3077
$ws.Range("I23").Value = 391.33334
$ws.Range("J23").Value = 781.3
$ws.Range("K23").Value = 1174.00002
$ws.Range("L23").Value = 2343.9
$ws.Range("M23").Value = -939.0000199999999
$ws.Range("N23").Value = -2813.9
$ws.Range("H37").Value = 116344.73
$ws.Range("J37").Value = 116344.73
$ws.Range("L37").Value = 349034.19
$ws.Range("N37").Value = -349258.19
$ws.Range("H38").Value = 61.07143
$ws.Range("I38").Value = 14
$ws.Range("K38").Value = 42
$ws.Range("M38").Value = 305
$ws.Range("H137").Value = 4727
$ws.Range("J137").Value = 4816.9165
$ws.Range("L137").Value = 14450.7495
$ws.Range("N137").Value = -24650.7495

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3940.0588
$ws.Range("I113").Value = 3498.818
$ws.Range("K113").Value = 3498.818
$ws.Range("M113").Value = -1328.818

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7062.375
$ws.Range("I7").Value = 5873.75
$ws.Range("J7").Value = 8251
$ws.Range("K7").Value = 5873.75
$ws.Range("L7").Value = 8251
$ws.Range("M7").Value = -5761.75
$ws.Range("N7").Value = -8475
$ws.Range("H40").Value = 4369.5
$ws.Range("J40").Value = 5000.3335
$ws.Range("L40").Value = 5000.3335
$ws.Range("N40").Value = -5272.3335
$ws.Range("H93").Value = 2521
$ws.Range("J93").Value = 2500
$ws.Range("L93").Value = 2500
$ws.Range("N93").Value = -4996
$ws.Range("H126").Value = 7062.375
$ws.Range("I126").Value = 5873.75
$ws.Range("J126").Value = 8251
$ws.Range("K126").Value = 17621.25
$ws.Range("L126").Value = 24753
$ws.Range("M126").Value = -15151.25
$ws.Range("N126").Value = -29693
$ws.Range("H133").Value = 49113.8
$ws.Range("J133").Value = 58325.668
$ws.Range("L133").Value = 58325.668
$ws.Range("N133").Value = -63385.668

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 401.5
$ws.Range("I3").Value = 401.5
$ws.Range("K3").Value = 401.5
$ws.Range("M3").Value = -287.5
$ws.Range("H96").Value = 6897
$ws.Range("I96").Value = 6252.7144
$ws.Range("K96").Value = 6252.7144
$ws.Range("M96").Value = -4879.7144
$ws.Range("H107").Value = 1022.38464
$ws.Range("I107").Value = 987.5294
$ws.Range("K107").Value = 2962.5882
$ws.Range("M107").Value = -1042.5882
$ws.Range("H132").Value = 1211.5238
$ws.Range("I132").Value = 1084.5883
$ws.Range("K132").Value = 3253.7649
$ws.Range("M132").Value = -723.7648999999997
$ws.Range("H136").Value = 11859.617
$ws.Range("I136").Value = 9778.412
$ws.Range("K136").Value = 29335.236
$ws.Range("M136").Value = -26785.236
